# Afry modell.xlsx - update long-term growth-rate assumptions
#
# Business change: in the "Modell" sheet, the projected growth applied to
# rows 4-8 (Omsetning/COGS buildup) for the forecast years in columns AD:AL
# is changed:
#   - column AD (2nd forecast year) now grows 5% off column AC (was 10%)
#   - columns AE:AL (3rd year onward) now grow 7.5% off the previous column
#     (was 10%)
# Column AC itself (first forecast year) is untouched.
#
# Every other formula in the workbook (SUM roll-ups, ratio rows, NPV, the
# Nøkkeltall sheet, ...) already references these cells, so recalculation
# propagates the new assumption automatically - we only need to touch the
# AD:AL formulas in rows 4-8.

$wb = $excel.ActiveWorkbook

$modell = $wb.Worksheets.Item("Modell")
$nokkeltall = $wb.Worksheets.Item("Nøkkeltall")

$rows = @(4, 5, 6, 7, 8)
foreach ($r in $rows) {
    # 2nd forecast year: was =AC*1.1, now =AC*1.05
    $modell.Range("AD${r}").Formula = "=AC${r}*1.05"
    # 3rd year through the last forecast year: was =prev*1.1, now =prev*1.075
    $modell.Range("AE${r}:AL${r}").Formula = "=AD${r}*1.075"
}

# Update the saved cursor/selection state to match the edited workbook.
$nokkeltall.Range("F19").Select()

$modell.Activate()
$modell.Range("AM31").Select()
